$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows to append to the bottom of the sheet (Indice 165 -> row 166,
# Indice 166 -> row 167), matching the style of the existing data rows.

# Reuse column A's "index" formatting (bold, bordered, centered) from the
# last existing data row instead of re-building it cell-by-cell, so it maps
# back onto the same cellXf the rest of the column already uses.
$ws.Range("A165").Copy()
$ws.Range("A166:A167").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @{
        Indice = 165
        pais = "romania"
        torneio = "liga-1"
        temporada = "2023-2024"
        data_partida = 45281.66666666666
        home = "FC Hermannstadt"
        home_ft_gols = 1
        away = "Sepsi Sf. Gheorghe"
        away_ft_gols = 1
        home_opening_odds = 2.38
        home_opening_data_hora = "17/12/2023 20:12"
        home_closing_odds = 2.62
        home_closing_data_hora = "21/12/2023 15:54"
        draw_opening_odds = 2.94
        draw_opening_data_hora = "17/12/2023 20:12"
        draw_closing_odds = 2.94
        draw_closing_data_hora = "21/12/2023 15:54"
        away_opening_odds = 3.19
        away_opening_data_hora = "17/12/2023 20:12"
        away_closing_odds = 3.11
        away_closing_data_hora = "21/12/2023 15:54"
        url_partida = "https://www.betexplorer.com/football/romania/liga-1/fc-hermannstadt-sepsi/fy1sRSy8/"
    },
    @{
        Indice = 166
        pais = "romania"
        torneio = "liga-1"
        temporada = "2023-2024"
        data_partida = 45281.78125
        home = "CFR Cluj"
        home_ft_gols = 4
        away = "U. Cluj"
        away_ft_gols = 0
        home_opening_odds = 1.59
        home_opening_data_hora = "17/12/2023 20:12"
        home_closing_odds = 1.97
        home_closing_data_hora = "21/12/2023 18:44"
        draw_opening_odds = 3.76
        draw_opening_data_hora = "17/12/2023 20:12"
        draw_closing_odds = 3.21
        draw_closing_data_hora = "21/12/2023 18:43"
        away_opening_odds = 5.3
        away_opening_data_hora = "17/12/2023 20:12"
        away_closing_odds = 4.39
        away_closing_data_hora = "21/12/2023 18:44"
        url_partida = "https://www.betexplorer.com/football/romania/liga-1/cfr-cluj-universitatea-cluj/b9RNZKu8/"
    }
)

$startRow = 166
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Indice

    $ws.Cells.Item($r, 2).Value = $row.pais
    $ws.Cells.Item($r, 3).Value = $row.torneio
    $ws.Cells.Item($r, 4).Value = $row.temporada

    $ws.Cells.Item($r, 5).Value = $row.data_partida
    $ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 6).Value = $row.home
    $ws.Cells.Item($r, 7).Value = $row.home_ft_gols
    $ws.Cells.Item($r, 8).Value = $row.away
    $ws.Cells.Item($r, 9).Value = $row.away_ft_gols

    $ws.Cells.Item($r, 10).Value = $row.home_opening_odds
    $ws.Cells.Item($r, 11).Value = $row.home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $row.home_closing_odds
    $ws.Cells.Item($r, 13).Value = $row.home_closing_data_hora

    $ws.Cells.Item($r, 14).Value = $row.draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $row.draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $row.draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $row.draw_closing_data_hora

    $ws.Cells.Item($r, 18).Value = $row.away_opening_odds
    $ws.Cells.Item($r, 19).Value = $row.away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $row.away_closing_odds
    $ws.Cells.Item($r, 21).Value = $row.away_closing_data_hora

    $ws.Cells.Item($r, 22).Value = $row.url_partida

    $r = $r + 1
}
